# Populate local copies of the Week 8 logs with Jesse's task-summary and
# activity-log-summary data (per commit: "copied over local copies of logs").

$wb = $excel.ActiveWorkbook

# ---- TASK SUMMARY SHEET ----
$ws3 = $wb.Worksheets.Item("TASK SUMMARY SHEET")

$ws3.Range("C1").Value = "Jesse Hare"
$ws3.Range("E1").Value = 8

$ws3.Range("A3").Value = "Project Build"
$ws3.Range("B3").Value = "Continue work on dynamic search/filter"
$ws3.Range("C3").Value = 6
$ws3.Range("D3").Value = 7
$ws3.Range("E3").Value = 0

$ws3.Range("A4").Value = "Project Build"
$ws3.Range("B4").Value = "Fix issues with the Search function"
$ws3.Range("C4").Value = 2
$ws3.Range("D4").Value = 4
$ws3.Range("E4").Value = 0

$ws3.Range("A5").Value = "Project Build"
$ws3.Range("B5").Value = "Implement new usability features and input validation"
$ws3.Range("C5").Value = 4
$ws3.Range("D5").Value = 5
$ws3.Range("E5").Value = 0

$ws3.Range("A6").Value = "Project Build"
$ws3.Range("B6").Value = "Testing of new features"
$ws3.Range("C6").Value = 4
$ws3.Range("D6").Value = 4
$ws3.Range("E6").Value = 0

# ---- ACTIVITY LOG SUMMARY SHEET ----
$ws4 = $wb.Worksheets.Item("ACTIVITY LOG SUMMARY SHEET")

$ws4.Range("D1").Value = "Jesse Hare"

$ws4.Range("A4").Value = "Project Build"
$ws4.Range("B4").Value = 15
$ws4.Range("C4").Value = 5

$wb.Application.Calculate()

# ---- Match the saved selection / active-sheet state ----
$ws4.Range("C6").Select() | Out-Null
$ws3.Activate() | Out-Null
$ws3.Range("C10").Select() | Out-Null
